$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "purpose" column (E) for rows 2-19 to "fullRNASEQ"
$ws.Range("E2:E19").Value = "fullRNASEQ"

# Reflect the final selection state from the edit
$ws.Range("D20:H26").Select()
